# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the "全部类型" sheet to the newly scraped values (gh-pages output
# regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1167
$ws1.Range("F3").Value  = 416
$ws1.Range("F4").Value  = 270
$ws1.Range("F6").Value  = 11
$ws1.Range("F7").Value  = 12276
$ws1.Range("F10").Value = 7
$ws1.Range("F11").Value = 138
$ws1.Range("F12").Value = 12055
$ws1.Range("F13").Value = 4810
$ws1.Range("F14").Value = 2814
$ws1.Range("F15").Value = 121
$ws1.Range("F16").Value = 52
$ws1.Range("F18").Value = 96
$ws1.Range("F19").Value = 947
$ws1.Range("F21").Value = 359
$ws1.Range("F22").Value = 165
$ws1.Range("F24").Value = 5215

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1167
$ws4.Range("F3").Value  = 416
$ws4.Range("F4").Value  = 270
$ws4.Range("F8").Value  = 11
$ws4.Range("F9").Value  = 12276
$ws4.Range("F12").Value = 7
$ws4.Range("F13").Value = 138
$ws4.Range("F14").Value = 12055
$ws4.Range("F15").Value = 4810
$ws4.Range("F16").Value = 2817
$ws4.Range("F17").Value = 121
$ws4.Range("F18").Value = 52
$ws4.Range("F20").Value = 96
$ws4.Range("F21").Value = 947
$ws4.Range("F23").Value = 359
$ws4.Range("F24").Value = 165
$ws4.Range("F26").Value = 5215
